$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Connector shapes (Google Shape ids 56,57,59,61,63,65,68) need their line
# weight bumped from 0.75pt (9525 EMU) to 1.5pt (19050 EMU) and the end
# arrowhead style switched from "stealth" to "triangle".
$connectorIds = @(56, 57, 59, 61, 63, 65, 68)
foreach ($id in $connectorIds) {
    $shapeIndex = $id - 53
    $cxn = $s.Shapes.Item($shapeIndex)
    $cxn.Line.Weight = 19050 / 12700
    $cxn.Line.EndArrowheadStyle = 2
}

# Shape 60 ("Cómputo de propiedades termodinámicas ...") shrinks and its
# text is trimmed.
$box = $s.Shapes.Item(60 - 53)
$box.Height = 831300 / 12700
$box.TextFrame.TextRange.Text = "Cómputo de propiedades termodinámicas"

# Shape 69 ("No" textbox) shifts right.
$noBox = $s.Shapes.Item(69 - 53)
$noBox.Left = 457.1575
